# Auto-generated edit script applying the Chocobo_Profits.xlsx diff
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H-N) for specific leves
# across the ALC, ARM, BSM, CUL, and LTW sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 34999.5
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 34999.5
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 34999.5
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -35223.5
$ws.Range("H10").Value = 29999
$ws.Range("J10").Value = 29999
$ws.Range("L10").Value = 29999
$ws.Range("N10").Value = -30585
$ws.Range("H14").Value = 34999.5
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 34999.5
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 34999.5
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -35381.5
$ws.Range("H112").Value = 1173.575
$ws.Range("I112").Value = 500
$ws.Range("J112").Value = 1209.0264
$ws.Range("K112").Value = 1500
$ws.Range("L112").Value = 3627.0792
$ws.Range("M112").Value = -392
$ws.Range("N112").Value = -5843.0792
$ws.Range("H139").Value = 33483.57
$ws.Range("J139").Value = 33483.57
$ws.Range("L139").Value = 33483.57
$ws.Range("N139").Value = -43763.57
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 833.0952
$ws.Range("I2").Value = 699.4666999999999
$ws.Range("J2").Value = 1167.1666
$ws.Range("K2").Value = 699.4666999999999
$ws.Range("L2").Value = 1167.1666
$ws.Range("M2").Value = -586.4666999999999
$ws.Range("N2").Value = -1393.1666
$ws.Range("H24").Value = 35000.332
$ws.Range("J24").Value = 35000.332
$ws.Range("L24").Value = 35000.332
$ws.Range("N24").Value = -35748.332
$ws.Range("H61").Value = 1662.5555
$ws.Range("I61").Value = 1662.5555
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1662.5555
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1450.5555
$ws.Range("N61").ClearContents()
$ws.Range("H100").Value = 35000.332
$ws.Range("J100").Value = 35000.332
$ws.Range("L100").Value = 35000.332
$ws.Range("N100").Value = -37164.332
$ws.Range("H110").Value = 1409.2
$ws.Range("I110").Value = 1364
$ws.Range("J110").Value = 1489.5555
$ws.Range("K110").Value = 1364
$ws.Range("L110").Value = 1489.5555
$ws.Range("M110").Value = 681
$ws.Range("N110").Value = -5579.5555
$ws.Range("H116").Value = 833.0952
$ws.Range("I116").Value = 699.4666999999999
$ws.Range("J116").Value = 1167.1666
$ws.Range("K116").Value = 699.4666999999999
$ws.Range("L116").Value = 1167.1666
$ws.Range("M116").Value = 1594.5333
$ws.Range("N116").Value = -5755.1666
$ws.Range("H132").Value = 3350.7856
$ws.Range("I132").Value = 1939
$ws.Range("J132").Value = 5233.1665
$ws.Range("K132").Value = 5817
$ws.Range("L132").Value = 15699.4995
$ws.Range("M132").Value = -3287
$ws.Range("N132").Value = -20759.4995
$ws.Range("H133").Value = 43332.668
$ws.Range("J133").Value = 43332.668
$ws.Range("L133").Value = 43332.668
$ws.Range("N133").Value = -48392.668
$ws.Range("H136").Value = 1662.5555
$ws.Range("I136").Value = 1662.5555
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4987.666499999999
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -2437.666499999999
$ws.Range("N136").ClearContents()
$ws.Range("H137").Value = 45565
$ws.Range("J137").Value = 45565
$ws.Range("L137").Value = 45565
$ws.Range("N137").Value = -55765
$ws.Range("H138").Value = 78990
$ws.Range("J138").Value = 78990
$ws.Range("L138").Value = 78990
$ws.Range("N138").Value = -89270
$ws.Range("H139").Value = 43999.137
$ws.Range("J139").Value = 43999.137
$ws.Range("L139").Value = 43999.137
$ws.Range("N139").Value = -54279.137
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 833.0952
$ws.Range("I3").Value = 699.4666999999999
$ws.Range("J3").Value = 1167.1666
$ws.Range("K3").Value = 699.4666999999999
$ws.Range("L3").Value = 1167.1666
$ws.Range("M3").Value = -585.4666999999999
$ws.Range("N3").Value = -1395.1666
$ws.Range("H59").Value = 48000
$ws.Range("J59").Value = 48000
$ws.Range("L59").Value = 48000
$ws.Range("N59").Value = -49694
$ws.Range("H94").Value = 799.3333
$ws.Range("I94").Value = 690.1818
$ws.Range("J94").Value = 2000
$ws.Range("K94").Value = 690.1818
$ws.Range("L94").Value = 2000
$ws.Range("M94").Value = -239.1818
$ws.Range("N94").Value = -2902
$ws.Range("H137").Value = 43666.668
$ws.Range("J137").Value = 43666.668
$ws.Range("L137").Value = 43666.668
$ws.Range("N137").Value = -53866.668
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1385.4762
$ws.Range("I5").Value = 474.0909
$ws.Range("J5").Value = 2388
$ws.Range("K5").Value = 1422.2727
$ws.Range("L5").Value = 7164
$ws.Range("M5").Value = -1310.2727
$ws.Range("N5").Value = -7388
$ws.Range("H38").Value = 147.25
$ws.Range("I38").Value = 57.5
$ws.Range("J38").Value = 165.2
$ws.Range("K38").Value = 172.5
$ws.Range("L38").Value = 495.6
$ws.Range("M38").Value = 174.5
$ws.Range("N38").Value = -1189.6
$ws.Range("H122").Value = 2396.64
$ws.Range("J122").Value = 2849.5264
$ws.Range("L122").Value = 25645.7376
$ws.Range("N122").Value = -30545.7376
$ws.Range("H131").Value = 786.4693600000001
$ws.Range("I131").Value = 442
$ws.Range("J131").Value = 804.9892599999999
$ws.Range("K131").Value = 1326
$ws.Range("L131").Value = 2414.96778
$ws.Range("M131").Value = 3714
$ws.Range("N131").Value = -12494.96778
$ws.Range("H135").Value = 1385.4762
$ws.Range("I135").Value = 474.0909
$ws.Range("J135").Value = 2388
$ws.Range("K135").Value = 4266.8181
$ws.Range("L135").Value = 21492
$ws.Range("M135").Value = -1731.8181
$ws.Range("N135").Value = -26562
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 230.38889
$ws.Range("I55").Value = 180.3
$ws.Range("J55").Value = 293
$ws.Range("K55").Value = 180.3
$ws.Range("L55").Value = 293
$ws.Range("M55").Value = -7.300000000000011
$ws.Range("N55").Value = -639
$ws.Range("H93").Value = 6947397
$ws.Range("I93").Value = 18521940
$ws.Range("J93").Value = 2671.2
$ws.Range("K93").Value = 18521940
$ws.Range("L93").Value = 2671.2
$ws.Range("M93").Value = -18520692
$ws.Range("N93").Value = -5167.2
$ws.Range("H132").Value = 5341
$ws.Range("I132").Value = 2200.7856
$ws.Range("J132").Value = 9337.637000000001
$ws.Range("K132").Value = 6602.3568
$ws.Range("L132").Value = 28012.911
$ws.Range("M132").Value = -4072.3568
$ws.Range("N132").Value = -33072.911
